$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date style used by the existing log rows (A2) so the new date
# cells (A9, A12) pick up the same number format / style index.
$ws.Range("A2").Copy()

# --- New task log rows appended below the existing entries (rows 2-7) ---
# Row 8
$ws.Range("B8").Value = "국립중앙도서관 API 발급"

# Row 9
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value = 45530
$ws.Range("B9").Value = "홈페이지 제작"

# Row 11 (entered before row 10 to match original authoring order)
$ws.Range("B11").Value = "책 디테일페이지 제작"

# Row 10
$ws.Range("B10").Value = "책 리스트페이지 제작"

# Row 12
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = 45531
$ws.Range("B12").Value = "작가 리스트페이지 제작"

# Row 13
$ws.Range("B13").Value = "작가 디테일페이지 제작"

# Row 14
$ws.Range("B14").Value = "세션 설정 및 저장"

# Row 15
$ws.Range("B15").Value = "사용자 계정 생성 및 권한설정"

# Row 17 (entered before row 16 to match original authoring order)
$ws.Range("B17").Value = "책대여갱신프로세스"

# Row 16
$ws.Range("B16").Value = "폼 생성"

# Row 18
$ws.Range("B18").Value = "unit testing 자동화"

# Row 20 - E20 entered before B20 to match original authoring order
$ws.Range("E20").Value = " "
$ws.Range("B20").Value = "국립중앙도서관 api 자료 다운로드"

# Row 21
$ws.Range("B21").Value = "오라클 DB에 저장"

# Row 22
$ws.Range("B22").Value = "오라클 - 로컬라이브러리 서버 연결"

# Match the saved view state (scroll position + active selection)
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("G14").Select()
